# Aggiornamento dati al 23 agosto 2021
# Appends rows 344-357 (dates 2021-08-10 .. 2021-08-23) to the data table,
# reusing the formatting of the last existing data row (343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44418, 3, 21, 186.0217911241031),
    @(44419, 2, 22, 194.8799716538223),
    @(44420, 2, 20, 177.1636105943839),
    @(44421, 4, 23, 203.7381521835415),
    @(44422, 2, 21, 186.0217911241031),
    @(44423, 0, 15, 132.8727079457879),
    @(44424, 1, 14, 124.0145274160687),
    @(44425, 1, 12, 106.2981663566303),
    @(44426, 0, 10, 88.58180529719195),
    @(44427, 1, 9, 79.72362476747276),
    @(44428, 1, 6, 53.14908317831517),
    @(44429, 3, 7, 62.00726370803437),
    @(44430, 2, 9, 79.72362476747276),
    @(44431, 0, 8, 70.86544423775356)
)

$lastRow = 343
$startRow = $lastRow + 1

# Copy the formatting of the last existing row onto the new rows in one shot.
$srcRow = $ws.Range("A$lastRow`:D$lastRow")
$endRow = $startRow + $data.Count - 1
$destRows = $ws.Range("A$startRow`:D$endRow")
$srcRow.Copy()
$destRows.PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
